$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph that follows the title.
$d.Paragraphs(2).Range.Delete()

# 2. Insert a new bold paragraph "Play Dazzle Me Megaways free - slot review 2021"
#    right before the final (image-prompt) paragraph.
$cnt = $d.Paragraphs.Count
$lastPara = $d.Paragraphs($cnt)
$lastPara.Range.InsertParagraphBefore()
$newPara = $d.Paragraphs($cnt)
$newParaRange = $d.Range($newPara.Range.Start, $newPara.Range.End)
$newParaRange.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Dazzle Me Megaways free - slot review 2021</w:t></w:r></w:p>')

# 3. Replace the old AI image-prompt text in the final paragraph with the
#    meta description text (keeping the existing italic run formatting).
$oldText = "Create a feature image for Dazzle Me Megaways that showcases the fun and energetic feel of the game. The image should be in a cartoon style and feature a happy Maya warrior with glasses to represent the adventurous and exciting aspect of the slot game. The warrior should be surrounded by colorful gemstones and fruits, which are symbols in the game. The background should have a bright and vibrant effect that brings out the lively experience of playing Dazzle Me Megaways."
$newText = "Discover the dazzling features, 99,225 ways to win and bonus rounds when you play Dazzle Me Megaways for free. Read our review on this slot game."
$d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
